$d = $word.ActiveDocument

# Fix "desenvolupado" -> "desarrollado"
$d.Content.Find.Execute("desenvolupado", $true, $false, $false, $false, $false,
                         $true, 1, $false, "desarrollado", 2)

# Merge the split "videotutorial" run by re-typing the whole sentence
$d.Content.Find.Execute("un videotutorial don", $true, $false, $false, $false, $false,
                         $true, 1, $false, "un videotutorial don", 2)
